$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)

# Simple numeric "想去人数" (interested-count) updates
$ws1.Range("F5").Value = 2273
$ws1.Range("F6").Value = 2
$ws1.Range("F8").Value = 73
$ws1.Range("F9").Value = 1671
$ws1.Range("F10").Value = 1671
$ws1.Range("F11").Value = 1381
$ws1.Range("F13").Value = 1430
$ws1.Range("F14").Value = 20
$ws1.Range("F16").Value = 630
$ws1.Range("F19").Value = 7374
$ws1.Range("F20").Value = 8217
$ws1.Range("F30").Value = 18
$ws1.Range("F34").Value = 1469
$ws1.Range("F35").Value = 251
$ws1.Range("F36").Value = 234
$ws1.Range("F44").Value = 255
$ws1.Range("F46").Value = 91
$ws1.Range("F47").Value = 189
$ws1.Range("F48").Value = 175
$ws1.Range("F49").Value = 12

# Row 31: a brand-new event is inserted here ("广播剧《恶人想要抢救一下》专场见面会"),
# pushing the previous occupants of rows 31-32 down by one, and the old cancelled
# row 33 event ("美漫超级英雄ONLY（取消）") drops out of the list entirely.
$ws1.Range("C31").Value = "北京·第19届IJOY漫展【广播剧《恶人想要抢救一下》专场见面会】"
$ws1.Range("E31").Value = "2024.10.02 11:00-10.02 15:30"
$ws1.Range("F31").Value = 0
$ws1.Range("G31").Value = 238
$ws1.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=92146"
$ws1.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202409/SFGR53ym1725853578783.jpeg"

# Row 32: now holds what was previously row 31 (IJOY 银发娘专场见面会); its date moves to 2024-10-02.
# Force the date-looking string to stay plain text (matches the rest of column B) instead of
# being auto-converted to a date serial by Excel, then drop the format override again.
$ws1.Range("B32").NumberFormat = "@"
$ws1.Range("B32").Value = "2024-10-02"
$ws1.Range("B32").ClearFormats()
$ws1.Range("C32").Value = "北京·第19届IJOY漫展【银发娘专场见面会】"
$ws1.Range("D32").Value = "天辰东路7号 北京国家会议中心"
$ws1.Range("E32").Value = "2024.10.02 12:25-10.02 16:30"
$ws1.Range("F32").Value = 9
$ws1.Range("G32").Value = 168
$ws1.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=91563"
$ws1.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202408/9Td79pPE1724928479521.jpeg"

# Row 33: now holds what was previously row 32 (明日方舟同人only-厮守序言)
$ws1.Range("C33").Value = "北京·明日方舟同人only-厮守序言"
$ws1.Range("D33").Value = "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园"
$ws1.Range("E33").Value = "2024.10.03 09:30-10.03 17:00"
$ws1.Range("F33").Value = 353
$ws1.Range("G33").Value = 60
$ws1.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=90959"
$ws1.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202408/rIGY6eyZ1723974119991.jpeg"

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 9
$ws2.Range("F5").Value = 64
$ws2.Range("F18").Value = 305

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 2636
$ws3.Range("F4").Value = 289
$ws3.Range("F6").Value = 18

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F9").Value = 2273
$ws4.Range("F11").Value = 73
$ws4.Range("F12").Value = 1671
$ws4.Range("F13").Value = 1671
$ws4.Range("F16").Value = 630
$ws4.Range("F17").Value = 9
$ws4.Range("F19").Value = 64
$ws4.Range("F21").Value = 7374
$ws4.Range("F22").Value = 8217
$ws4.Range("F27").Value = 18
$ws4.Range("F28").Value = 9
$ws4.Range("F29").Value = 251
$ws4.Range("F30").Value = 234
$ws4.Range("F44").Value = 255
$ws4.Range("F46").Value = 91
$ws4.Range("F47").Value = 189
$ws4.Range("F49").Value = 305

